# Applies the "Updated symbol list" GitHub Actions edit to the crypto price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price / Volume(1h) columns store numeric- and percentage-looking values as
# plain text (inline strings) in the source file. Force a text number format on
# those ranges before writing the new values so Excel does not reinterpret them
# as numbers/percentages. (Two separate contiguous ranges are used instead of a
# single multi-area union, since a union Range only applies property writes to
# its first area.)
$textRange1 = $ws.Range("D2:E27")
$textRange2 = $ws.Range("D39:E51")
$textRange1.NumberFormat = "@"
$textRange2.NumberFormat = "@"

# --- Coin / Link columns ---
# GateToken (previously listed at row 17) moved up to row 8, shifting the other
# exchange-token rows (MXToken .. LEO) down by one position.
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Price / Volume(1h) columns ---
$ws.Range("D2").Value = "304.02"
$ws.Range("E2").Value = "2.79%"
$ws.Range("D3").Value = "35.09"
$ws.Range("E3").Value = "12.78%"
$ws.Range("D4").Value = "5.181"
$ws.Range("E4").Value = "4.57%"
$ws.Range("D5").Value = "0.07797"
$ws.Range("E5").Value = "5.08%"
$ws.Range("D6").Value = "2.365"
$ws.Range("E6").Value = "10.56%"
$ws.Range("D7").Value = "8.037"
$ws.Range("E7").Value = "3.63%"
$ws.Range("D8").Value = "3.950"
$ws.Range("E8").Value = "5.37%"
$ws.Range("D9").Value = "0.9359"
$ws.Range("E9").Value = "2.33%"
$ws.Range("D10").Value = "0.09843"
$ws.Range("E10").Value = "12.13%"
$ws.Range("D11").Value = "0.1790"
$ws.Range("E11").Value = "5.05%"
$ws.Range("D12").Value = "0.08590"
$ws.Range("E12").Value = "3.64%"
$ws.Range("D13").Value = "0.03316"
$ws.Range("E13").Value = "5.33%"
$ws.Range("D14").Value = "0.09914"
$ws.Range("E14").Value = "-1.70%"
$ws.Range("D15").Value = "0.001502"
$ws.Range("E15").Value = "-0.27%"
$ws.Range("D16").Value = "0.005853"
$ws.Range("E16").Value = "0.99%"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").Value = "-1.34%"
$ws.Range("D18").Value = "2.189"
$ws.Range("E18").Value = "5.34%"
$ws.Range("E19").Value = "1.12%"
$ws.Range("D20").Value = "0.1295"
$ws.Range("E20").Value = "0.16%"
$ws.Range("D21").Value = "4.349"
$ws.Range("E21").Value = "9.15%"
$ws.Range("D22").Value = "0.2301"
$ws.Range("E22").Value = "8.01%"
$ws.Range("D23").Value = "0.04623"
$ws.Range("E23").Value = "1.63%"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").Value = "0.00%"
$ws.Range("D25").Value = "0.004384"
$ws.Range("E25").Value = "-5.29%"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "-0.04%"
$ws.Range("D27").Value = "0.0003395"
$ws.Range("D39").Value = "0.01804"
$ws.Range("E39").Value = "11.38%"
$ws.Range("D40").Value = "0.04816"
$ws.Range("E40").Value = "7.34%"
$ws.Range("D41").Value = "0.007798"
$ws.Range("E41").Value = "7.11%"
$ws.Range("D42").Value = "0.1415"
$ws.Range("E42").Value = "6.14%"
$ws.Range("D43").Value = "0.009054"
$ws.Range("E43").Value = "0.86%"
$ws.Range("D44").Value = "0.002086"
$ws.Range("E44").Value = "6.25%"
$ws.Range("D45").Value = "0.009416"
$ws.Range("E45").Value = "0.26%"
$ws.Range("D46").Value = "0.00006110"
$ws.Range("E46").Value = "0.28%"
$ws.Range("E47").Value = "0.02%"
$ws.Range("D48").Value = "2.861"
$ws.Range("E48").Value = "27.66%"
$ws.Range("D49").Value = "0.002001"
$ws.Range("E49").Value = "-30.97%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "0.02%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "0.02%"

# Remove the temporary text style again so the cells fall back to the default
# (unstyled) cell format, matching the original workbook.
$textRange1.Style = "Normal"
$textRange2.Style = "Normal"
